# Division Master.xlsx — header rename (FullName -> Name) and the
# cosmetic fallout of re-saving the workbook afterwards: the A/B column
# widths re-fit to the shorter header text and the selection moved on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "FullName" -> "Name" header text in A1
$ws.Range("A1").Value = "Name"

# Columns re-fit to the new header text
$ws.Columns.Item(1).ColumnWidth = 5.5
$ws.Columns.Item(2).ColumnWidth = 10

# Active selection left on D5
[void]$ws.Range("D5").Select()
